# Update the "Förändrad" date column (C) for all data rows (2-43)
# from 45747 (2025-03-31) to 45749 (2025-04-02).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

for ($row = 2; $row -le 43; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45747) {
        $cell.Value2 = 45749
    }
}
